$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.658546805381775
$ws.Range("B1").Value = 3.486228466033936
$ws.Range("C1").Value = 4.056802272796631
$ws.Range("D1").Value = 1.27882182598114
$ws.Range("E1").Value = 0.7499282360076904
